# Replace the "has_grimoire" label with "zaubern" (new method to create grimoire)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "zaubern"

# Update the selected cell/cursor position as recorded in the saved file
$ws.Range("B30").Select()
